# "fixed experiment pipeline and did some code structure changes"
#
# - Drop the old "Experiment_2" worksheet (its results are superseded).
# - Rename "Experiment_1" to "Experiment_2" (it becomes the sole experiment sheet).
# - Replace its results table: drop the alpha/beta sweep columns and keep
#   only Similarity / Inference_Time, with refreshed values.

$wb = $excel.ActiveWorkbook

# Remove the stale "Experiment_2" sheet entirely.
$wb.Worksheets.Item("Experiment_2").Delete()

# The former "Experiment_1" sheet becomes the new "Experiment_2".
$ws = $wb.Worksheets.Item("Experiment_1")
$ws.Name = "Experiment_2"

# Drop the old alpha/beta columns - only Similarity/Inference_Time remain.
$ws.Range("C1:D2").Clear()

# Refresh the header row.
$ws.Range("A1").Value = "Similarity"
$ws.Range("B1").Value = "Inference_Time"

# Refresh the result row.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 837.4010593891144
